# Apply crypto price/volume updates per the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.218.92"
$ws.Range("E2").Value = "  +5.57%  "

$ws.Range("D3").Value = "1.909.75"
$ws.Range("E3").Value = "  +1.97%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  -0.56%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "329.15"
$ws.Range("E5").Value = "  +4.62%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  -0.59%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5178"
$ws.Range("E7").Value = "  +1.97%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4043"
$ws.Range("E8").Value = "  +3.57%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08475"
$ws.Range("E9").Value = "  +1.30%  "

$ws.Range("B10").Value = "OKB"
$ws.Range("C10").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.83"
$ws.Range("E10").Value = "  +1.31%  "

$ws.Range("B11").Value = "Polygon"
$ws.Range("C11").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.123"
$ws.Range("E11").Value = "  +1.63%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "23.24"
$ws.Range("E12").Value = "  +14.43%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.415"
$ws.Range("E13").Value = "  +3.74%  "

$ws.Range("E14").Value = "  +2.31%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.371"
$ws.Range("E15").Value = "  +1.77%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.003"
$ws.Range("E16").Value = "  -0.61%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "94.85"
$ws.Range("E17").Value = "  +1.81%  "

$ws.Range("E18").Value = "  +1.27%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06696"
$ws.Range("E19").Value = "  -0.16%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.40"
$ws.Range("E20").Value = "  +4.52%  "

$ws.Range("E21").Value = "  -0.65%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.005"
$ws.Range("E22").Value = "  +1.36%  "

$ws.Range("D23").Value = "30.220.12"
$ws.Range("E23").Value = "  +5.52%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.34"
$ws.Range("E24").Value = "  +2.25%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.225"
$ws.Range("E25").Value = "  +0.26%  "

$ws.Range("D26").Value = "2.135.46"
$ws.Range("E26").Value = "  +2.33%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.03"
$ws.Range("E27").Value = "  +7.11%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "162.30"
$ws.Range("E28").Value = "  +2.96%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.387"
$ws.Range("E29").Value = "  -0.72%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "129.80"
$ws.Range("E30").Value = "  +2.92%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.098"
$ws.Range("E31").Value = "  +5.36%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1063"
$ws.Range("E32").Value = "  +2.76%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.988"
$ws.Range("E33").Value = "  +3.67%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.647"
$ws.Range("E34").Value = "  +0.58%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02484"
$ws.Range("E35").Value = "  +1.50%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06557"
$ws.Range("E36").Value = "  +0.55%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2203"
$ws.Range("E37").Value = "  +2.02%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.187"
$ws.Range("E38").Value = "  +2.95%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.228"
$ws.Range("E39").Value = "  +2.97%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "11.98"
$ws.Range("E40").Value = "  +7.79%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.797"
$ws.Range("E41").Value = "  -2.13%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6520"
$ws.Range("E42").Value = "  +2.37%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.236"
$ws.Range("E43").Value = "  +0.10%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6136"
$ws.Range("E44").Value = "  +2.56%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.34"
$ws.Range("E45").Value = "  +2.43%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.740"
$ws.Range("E46").Value = "  +1.77%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.065"
$ws.Range("E47").Value = "  +3.27%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.244"
$ws.Range("E48").Value = "  +1.80%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "124.90"
$ws.Range("E49").Value = "  +2.55%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.161"
$ws.Range("E50").Value = "  -3.54%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "79.33"
$ws.Range("E51").Value = "  +4.29%  "
